$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 18; this shifts the existing rows 18 and 19
# down to rows 19 and 20 respectively, preserving all of their data/formatting.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new weekly price entry.
$ws.Cells.Item(18, 1).Value = 11
$ws.Cells.Item(18, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(18, 3).Value = "Bíobío"
$ws.Cells.Item(18, 4).Value = 45265
$ws.Cells.Item(18, 5).Value = 8
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100103
$ws.Cells.Item(18, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(18, 9).Value = 100103003
$ws.Cells.Item(18, 10).Value = "Damasco"
$ws.Cells.Item(18, 11).Value = "Dina"
$ws.Cells.Item(18, 12).Value = "Primera"
$ws.Cells.Item(18, 13).Value = 150
$ws.Cells.Item(18, 14).Value = 25000
$ws.Cells.Item(18, 15).Value = 25000
$ws.Cells.Item(18, 16).Value = 25000
$ws.Cells.Item(18, 17).Value = "`$/caja 10 kilos"
$ws.Cells.Item(18, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(18, 19).Value = 2500
$ws.Cells.Item(18, 20).Value = 10

# Match the date-cell number format used by the other rows in column D.
$ws.Cells.Item(18, 4).NumberFormat = $ws.Cells.Item(19, 4).NumberFormat

$wb.Save()
